# "thêm đơn hàng trực tiếp vào kpi" — add "direct order" KPI columns
# alongside the existing "indirect order" KPI columns on the "KPI san pham" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPI san pham")

# 1) Relabel the existing KPI metric headers (D4:G4) to clarify they are the
#    "indirect order" (đơn gián tiếp) figures.
$ws.Range("D4").Value = "Sản lượng theo đơn gián tiếp"
$ws.Range("E4").Value = "Doanh thu theo đơn gián tiếp"
$ws.Range("F4").Value = "Số đơn hàng theo đơn gián tiếp"
$ws.Range("G4").Value = "Số đại lý theo đơn gián tiếp"

# 2) Add four new columns (H:K) for the "direct order" (đơn trực tiếp) figures,
#    mirroring the formatting of the existing metric columns.
$ws.Range("G4").Copy() | Out-Null
$ws.Range("H4:K4").PasteSpecial(-4122) | Out-Null
$ws.Range("H4").Value = "Sản lượng theo đơn trực tiếp"
$ws.Range("I4").Value = "Doanh thu theo đơn trực tiếp"
$ws.Range("J4").Value = "Số đơn hàng theo đơn trực tiếp"
$ws.Range("K4").Value = "Số đại lý theo đơn trực tiếp"

$ws.Range("G5").Copy() | Out-Null
$ws.Range("H5:K5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# 3) Size the new columns to fit their (longer) header text (best-fit-style
#    widths, matching the values Excel's own AutoFit would have produced).
$ws.Columns.Item(8).ColumnWidth = 26.8
$ws.Columns.Item(9).ColumnWidth = 26.8
$ws.Columns.Item(10).ColumnWidth = 29
$ws.Columns.Item(11).ColumnWidth = 25.2

# 4) Match the author's final selection.
$ws.Range("D18").Select() | Out-Null
